$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, shifting existing rows 107:132 down to 108:133
$ws.Rows("107").Insert()

# Populate the new row 107 with a fresh data record (copy of the row pattern
# with new Fecha / Precio / Origen / Precio $/Kg values)
$ws.Range("A107").Value = 5
$ws.Range("B107").Value = "Macroferia Regional de Talca"
$ws.Range("C107").Value = "Maule"
$ws.Range("D107").Value = 45204
$ws.Range("E107").Value = 7
$ws.Range("F107").Value = 100112040
$ws.Range("G107").Value = "Cilantro"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 9000
$ws.Range("L107").Value = 9000
$ws.Range("M107").Value = 9000
$ws.Range("N107").Value = "$/caja 36 atados"
$ws.Range("O107").Value = "Región Metropolitana"
$ws.Range("P107").Value = 250
$ws.Range("Q107").Value = 36
$ws.Range("R107").Value = "Hortaliza"
